$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# Each hashtable: target cell (1-based Row/Col) plus the four lattice-exercise text lines
# (the constant "  ----" divider line is appended in the loop below).
$cellsData = @(
    @{Row=1; Col=1; Top="22 x 64"; Digits="  6    4"; B1="2|    |"; B2="2|    |"}
    @{Row=1; Col=2; Top="94 x 36"; Digits="  3    6"; B1="9|    |"; B2="4|    |"}
    @{Row=1; Col=3; Top="57 x 42"; Digits="  4    2"; B1="5|    |"; B2="7|    |"}
    @{Row=2; Col=1; Top="52 x 50"; Digits="  5    0"; B1="5|    |"; B2="2|    |"}
    @{Row=2; Col=2; Top="31 x 55"; Digits="  5    5"; B1="3|    |"; B2="1|    |"}
    @{Row=2; Col=3; Top="66 x 54"; Digits="  5    4"; B1="6|    |"; B2="6|    |"}
    @{Row=3; Col=1; Top="95 x 83"; Digits="  8    3"; B1="9|    |"; B2="5|    |"}
    @{Row=3; Col=2; Top="96 x 95"; Digits="  9    5"; B1="9|    |"; B2="6|    |"}
    @{Row=3; Col=3; Top="52 x 53"; Digits="  5    3"; B1="5|    |"; B2="2|    |"}
    @{Row=4; Col=1; Top="84 x 24"; Digits="  2    4"; B1="8|    |"; B2="4|    |"}
    @{Row=4; Col=2; Top="30 x 95"; Digits="  9    5"; B1="3|    |"; B2="0|    |"}
    @{Row=4; Col=3; Top="59 x 99"; Digits="  9    9"; B1="5|    |"; B2="9|    |"}
    @{Row=5; Col=1; Top="11 x 61"; Digits="  6    1"; B1="1|    |"; B2="1|    |"}
    @{Row=5; Col=2; Top="44 x 57"; Digits="  5    7"; B1="4|    |"; B2="4|    |"}
    @{Row=5; Col=3; Top="36 x 30"; Digits="  3    0"; B1="3|    |"; B2="6|    |"}
)

foreach ($c in $cellsData) {
    $newText = $c.Top + $nl + $c.Digits + $nl + "  ----" + $nl + $c.B1 + $nl + $c.B2
    $cell = $t.Cell($c.Row, $c.Col)
    $cell.Range.Text = $newText
}

Write-Output "Updated lattice multiplication table cells"
